# Update Sheets via scheduled runner
# Refreshes market price / profit figures for several leve rows across sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1250513.5
$ws.Range("I38").Value = 81.59999999999999
$ws.Range("J38").Value = 3334566.8
$ws.Range("K38").Value = 244.8
$ws.Range("L38").Value = 10003700.4
$ws.Range("M38").Value = 127.2
$ws.Range("N38").Value = -10004444.4

$ws.Range("H41").Value = 1183.2354
$ws.Range("I41").Value = 1406.3334
$ws.Range("J41").Value = 647.8
$ws.Range("K41").Value = 1406.3334
$ws.Range("L41").Value = 647.8
$ws.Range("M41").Value = -966.3334
$ws.Range("N41").Value = -1527.8

$ws.Range("H53").Value = 2126.818
$ws.Range("I53").Value = 4540
$ws.Range("K53").Value = 4540
$ws.Range("M53").Value = -3903

$ws.Range("H58").Value = 617.2727
$ws.Range("I58").Value = 465
$ws.Range("J58").Value = 800
$ws.Range("K58").Value = 1395
$ws.Range("L58").Value = 2400
$ws.Range("M58").Value = -1245
$ws.Range("N58").Value = -2700

$ws.Range("H64").Value = 4759.5
$ws.Range("I64").Value = 3409.3157
$ws.Range("K64").Value = 3409.3157
$ws.Range("M64").Value = -3161.3157

$ws.Range("H67").Value = 4759.5
$ws.Range("I67").Value = 3409.3157
$ws.Range("K67").Value = 3409.3157
$ws.Range("M67").Value = -2551.3157

$ws.Range("H74").Value = 3424.95
$ws.Range("I74").Value = 3285.5715
$ws.Range("K74").Value = 3285.5715
$ws.Range("M74").Value = -2349.5715

$ws.Range("H77").Value = 3424.95
$ws.Range("I77").Value = 3285.5715
$ws.Range("K77").Value = 16427.8575
$ws.Range("M77").Value = -11747.8575

$ws.Range("H111").Value = 1870
$ws.Range("J111").Value = 1494.4445
$ws.Range("L111").Value = 4483.333500000001
$ws.Range("N111").Value = -10617.3335

$ws.Range("H137").Value = 14288051
$ws.Range("I137").Value = 1371.4286
$ws.Range("J137").Value = 28574730
$ws.Range("K137").Value = 4114.2858
$ws.Range("L137").Value = 85724190
$ws.Range("M137").Value = -1564.2858
$ws.Range("N137").Value = -85729290

$ws.Range("H138").Value = 4631466.5
$ws.Range("I138").Value = 1363.174
$ws.Range("J138").Value = 8066704.5
$ws.Range("K138").Value = 4089.522
$ws.Range("L138").Value = 24200113.5
$ws.Range("M138").Value = 1050.478
$ws.Range("N138").Value = -24210393.5

$ws.Range("H139").Value = 31475
$ws.Range("J139").Value = 31475
$ws.Range("L139").Value = 31475
$ws.Range("N139").Value = -41755

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 500
$ws.Range("I6").Value = 500
$ws.Range("K6").Value = 500
$ws.Range("M6").Value = -327

$ws.Range("H61").Value = 38466524
$ws.Range("I61").Value = 62505224
$ws.Range("J61").Value = 4602.8
$ws.Range("K61").Value = 62505224
$ws.Range("L61").Value = 4602.8
$ws.Range("M61").Value = -62505012
$ws.Range("N61").Value = -5026.8

$ws.Range("H132").Value = 25003838
$ws.Range("I132").Value = 50002904
$ws.Range("K132").Value = 150008712
$ws.Range("M132").Value = -150006182

$ws.Range("H136").Value = 38466524
$ws.Range("I136").Value = 62505224
$ws.Range("J136").Value = 4602.8
$ws.Range("K136").Value = 187515672
$ws.Range("L136").Value = 13808.4
$ws.Range("M136").Value = -187513122
$ws.Range("N136").Value = -18908.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4442
$ws.Range("I105").Value = 3341.125
$ws.Range("J105").Value = 4780.731
$ws.Range("K105").Value = 3341.125
$ws.Range("L105").Value = 4780.731
$ws.Range("M105").Value = -1594.125
$ws.Range("N105").Value = -8274.731

$ws.Range("H107").Value = 1462.5807
$ws.Range("I107").Value = 1524
$ws.Range("K107").Value = 1524
$ws.Range("M107").Value = 396

$ws.Range("H134").Value = 4117.136
$ws.Range("I134").Value = 2392.4
$ws.Range("J134").Value = 5554.4165
$ws.Range("K134").Value = 7177.200000000001
$ws.Range("L134").Value = 16663.2495
$ws.Range("M134").Value = -4642.200000000001
$ws.Range("N134").Value = -21733.2495

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 17552334
$ws.Range("I31").Value = 11846.25
$ws.Range("J31").Value = 47621740
$ws.Range("K31").Value = 11846.25
$ws.Range("L31").Value = 47621740
$ws.Range("M31").Value = -11551.25
$ws.Range("N31").Value = -47622330

$ws.Range("H34").Value = 17552334
$ws.Range("I34").Value = 11846.25
$ws.Range("J34").Value = 47621740
$ws.Range("K34").Value = 11846.25
$ws.Range("L34").Value = 47621740
$ws.Range("M34").Value = -11644.25
$ws.Range("N34").Value = -47622144

$ws.Range("H58").Value = 3658.4285
$ws.Range("I58").Value = 1857
$ws.Range("J58").Value = 5459.857
$ws.Range("K58").Value = 1857
$ws.Range("L58").Value = 5459.857
$ws.Range("M58").Value = -1654
$ws.Range("N58").Value = -5865.857

$ws.Range("H122").Value = 1489.2307
$ws.Range("J122").Value = 1875.4
$ws.Range("L122").Value = 5626.200000000001
$ws.Range("N122").Value = -10526.2

$ws.Range("H136").Value = 3658.4285
$ws.Range("I136").Value = 1857
$ws.Range("J136").Value = 5459.857
$ws.Range("K136").Value = 5571
$ws.Range("L136").Value = 16379.571
$ws.Range("M136").Value = -3021
$ws.Range("N136").Value = -21479.571

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 209
$ws.Range("I2").Value = 336.66666
$ws.Range("J2").Value = 17.5
$ws.Range("K2").Value = 2019.99996
$ws.Range("L2").Value = 105
$ws.Range("M2").Value = -1906.99996
$ws.Range("N2").Value = -331

$ws.Range("H23").Value = 82.38461
$ws.Range("I23").Value = 58.875
$ws.Range("K23").Value = 176.625
$ws.Range("M23").Value = 58.375

$ws.Range("H92").Value = 2000
$ws.Range("J92").Value = 2000
$ws.Range("L92").Value = 6000
$ws.Range("N92").Value = -8496

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 2500
$ws.Range("J38").Value = 2500
$ws.Range("L38").Value = 2500
$ws.Range("N38").Value = -3426

$ws.Range("H48").Value = 230000
$ws.Range("I48").Value = 230000
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 230000
$ws.Range("L48").Value = 0
$ws.Range("M48").Value = -229515
$ws.Range("N48").ClearContents()

$ws.Range("H49").Value = 39800
$ws.Range("J49").Value = 39800
$ws.Range("L49").Value = 39800
$ws.Range("N49").Value = -40168

$ws.Range("H132").Value = 6120.2104
$ws.Range("I132").Value = 5514.1665
$ws.Range("J132").Value = 6399.923
$ws.Range("K132").Value = 16542.4995
$ws.Range("L132").Value = 19199.769
$ws.Range("M132").Value = -14012.4995
$ws.Range("N132").Value = -24259.769

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H41").Value = 186666.67
$ws.Range("I41").Value = 186666.67
$ws.Range("K41").Value = 186666.67
$ws.Range("M41").Value = -186228.67

$ws.Range("H139").Value = 39810.445
$ws.Range("J139").Value = 39705.5
$ws.Range("L139").Value = 39705.5
$ws.Range("N139").Value = -49985.5
